$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: mint a brand-new numbering definition (abstractNum + num) by
# applying a list template to a throw-away paragraph at the end of the
# document, tuning the per-level formats to match a classic
# decimal / lowerLetter / lowerRoman hybrid-multilevel list, then deleting
# the scratch paragraph again (the numbering definition itself survives in
# numbering.xml even once nothing in the body still references it... but we
# will reference it again in step 2 with the same numId).
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertBefore("TEMP-LIST-SEED")
$seedPara = $d.Paragraphs.Last
$seedPara.Style = "List Paragraph"

$listGallery = $word.ListGalleries.Item(2)
$listTemplate = $listGallery.ListTemplates.Item(1)
$seedPara.Range.ListFormat.ApplyListTemplate($listTemplate)

$appliedTemplate = $seedPara.Range.ListFormat.ListTemplate
$levelFormats = @("decimal", "lowerLetter", "lowerRoman", "decimal", "lowerLetter", "lowerRoman", "decimal", "lowerLetter", "lowerRoman")
$levelStyles  = @(0, 4, 2, 0, 4, 2, 0, 4, 2)
for ($i = 1; $i -le $appliedTemplate.ListLevels.Count; $i++) {
    $lvl = $appliedTemplate.ListLevels.Item($i)
    $lvl.NumberStyle = $levelStyles[$i - 1]
    $lvl.NumberFormat = "%$i."
}

# remove the scratch paragraph; the minted numbering definition remains
$seedPara.Range.Delete()

# ---------------------------------------------------------------------------
# Step 2: insert the real "Firebase Setup" content as a raw OOXML fragment
# right after the last paragraph of the document (before the sectPr), via
# InsertXML so every run/paragraph property lands exactly as authored.
# ---------------------------------------------------------------------------
$insertionPoint = $d.Content
$insertionPoint.Collapse(0)

$bodyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:ind w:left="360"/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:ind w:left="360"/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:ind w:left="360"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:lastRenderedPageBreak/>
              <w:t>Firebase Setup:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="5"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Create new Project in firebase and create a new web app in the project</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="5"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Navigate to Project Settings and scroll down to </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>SDK setup and configuration</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="5"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Click on </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>C</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t xml:space="preserve">onfig </w:t>
            </w:r>
            <w:r>
              <w:t>and copy the configuration object shown</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="5"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Create a new file </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>firebase.js</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>in project folder and paste it</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="5"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Install firebase in the project folder using the commands</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t>For the first time use this to install firebase globally:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>(</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>npm install -g firebase-tools)</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Next, In the project directory, </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>npm i firebase</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t>Check this link for user register and sign in process using firebase</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t>https://firebase.google.com/docs/auth/web/start</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:ind w:left="360"/>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($bodyXml)
